# Automatic update of files.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update "Förändrad" (changed) date column C for rows 2-11 from 45208 to 45212
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45212
}

# Update hyperlink formulas in row 2 (A 33491-2023) with corrected filenames
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2104/artfynd/A 33491-2023 artfynd.xlsx", "A 33491-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2104/kartor/A 33491-2023 karta.png", "A 33491-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2104/klagomål/A 33491-2023 fsc-klagomål.docx", "A 33491-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2104/klagomålsmail/A 33491-2023 fsc-klagomål mail.docx", "A 33491-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2104/tillsyn/A 33491-2023 tillsynsbegäran.docx", "A 33491-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2104/ti,llsynsmail/A 33491-2023 tillsynsbegäran mail.docx", "A 33491-2023")'
